$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text updates (shared strings) ---

# D7: intro text for "Spell Crafting" quest row - prepend a "Note" paragraph
$ws.Range("D7").Value = @'
<p><strong>Note: </strong>This quest can take up to 2 hours to complete. Do not feel like you have to rush through it.</p><p>Lets start learning about spells.</p><p>You have just spent a ton f time crafting weapons and Armour, we will repeat this to craft spells. Spells are useful for caster classes, as well as any one who wants to use them.</p><p>Damage spells and Staves (Two Handed weapons) can raise the characters intelligence, both of which can be bought from the shop or crafted.</p><p>Healing spells are great for characters who want to do Cast and attack or Attack and Cast, like Prophets.</p><p>To get the quest item required, you will need to kill: Umbering Spirit Lord on Surface. This creature is further down the list and may require you to upgrade your gear through the shop before being able to take him down. This creature has a 15% chance to drop the item, so exploration might be a good choice here.</p>
'@

# C8: intro_text for "Go To Labyrinth" quest row - add Enchanting teaser sentence
$ws.Range("C8").Value = @'
You have been under the Enchantress’s tutelage for weeks now. She has been teaching the you about the art of magic and how it can be used to weave together powerful spells.<br /> <br /> She has told you that if you wait for her in the forest outside of town she will teach you the next step: Enchanting. She explained it is the use of powerful magic that can imbue items, such as weapons, armour, rings and even other spells with potent abilities that can turn the tide of a battle, similar to the items you find on the enemies corpses.<br /> <br /> You are waiting for her in the forest outside of the town you have been practically living in. She doesn’t appear and some time goes by.<br /> <br /> The old man appears, The Guide. Almost from thin air, one moment nothing, next he’s in front of you.<br /> <br /> “She isn’t coming.” He states.<br /> <br /> “Something foul is afoot and I need you to do some investigating for me. Find the Key of labyrinth child. Do it quickly.”
'@

# D8: instructions for "Go To Labyrinth" quest row - fix typo/wording, add "on Surface"
$ws.Range("D8").Value = @'
<p>Time to stop training and fight a monster. This will drop weather you use exploration or not. Find the Key of Labyrinth by fighting the Labyrinth Fiend on Surface</p><p>To make it easier lets also raise our looting skill.</p><p>Once you have the key, lets:</p><p><strong>Desktop:</strong></p><p><strong>- </strong>Click traverse under the map.</p><p>- Select Labyrinth from the plane drop down.</p><p>- Click traverse.</p><p>- Welcome to labyrinth. Monsters down here have the same strength as surface but are new.</p><p>- Kill enough until your faction level with surface is level 1.</p><p><strong>Mobile:</strong></p><p>- From the action drop down select Map Movement</p><p><strong>- </strong>Click traverse under the map.</p><p>- Select Labyrinth from the plane drop down.</p><p>- Click traverse.</p><p>- Welcome to labyrinth. Monsters down here have the same strength as surface but are new.</p><p>- Kill enough until your faction level with surface is level 1.</p>
'@

# --- Numeric updates ---

# "Go To Labyrinth" row (id 10, row 8): required_skill_level 250 -> 25, required_secondary_skill_level 50 -> 10
$ws.Range("G8").Value = 25
$ws.Range("I8").Value = 10

# "Enchanting is key" row (id 11, row 9): required_skill_level 20 -> 15, required_secondary_skill_level 20 -> 12
$ws.Range("G9").Value = 15
$ws.Range("I9").Value = 12
